# The "Date: " paragraph currently reads "2023-10-8" and needs to become
# "2023-10-08", with the corrected text split across three runs:
#   "2023-10-", "0", "8"
# (matching the way Word itself would split a run when a character is
# typed in the middle of existing text).

$d = $word.ActiveDocument

# Step 1: find the existing date text and shrink it down to just its
# leading portion, "2023-10-". This rewrites the original single run
# "2023-10-8" so that it now only contains "2023-10-".
$dateRange = $d.Content
$dateRange.Find.Execute("2023-10-8")
$dateRange.Text = "2023-10-"

# Step 2: insert the remaining two characters as their own separate runs
# right after it, one at a time, so the final text reads "2023-10-08"
# split across three runs: "2023-10-", "0", "8".
$insertPos = $dateRange.End

$zero = $d.Range($insertPos, $insertPos)
$zero.InsertAfter("0")

$eight = $d.Range($insertPos + 1, $insertPos + 1)
$eight.InsertAfter("8")
